# Update GDB Results Tool citation
#
# Row 79 (cell A79) holds the short-form source citation for the
# mortality data; Row 80 (cell A80) holds the longer methodology note.
# This edit expands the citation in A79 with the full GBD Results Tool
# reference details (publisher, year, access date, URL), matching the
# author's commit "Update GDB Results Tool citation". A80 keeps the
# original methodology note text unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Demographics")

$ws.Range("A79").Value = "UN World Population Prospects 2019, File MORT/17-2: Abridged life tables, for males/females, 1950-2100; GHDx GBD Results Tool: Global Burden of Disease Collaborative Network, Global Burden of Disease Study 2017 (GBD 2017) Results. Seattle, United States: Institute for Health Metrics and Evaluation (IHME), 2018. Available from http://ghdx.healthdata.org/gbd-results-tool. Accessed Feb 21, 2020."

$ws.Range("A80").Value = "Used UN values from 1950-1990 before the widespread HIV epidemic. Used the GHDx GBD Results Tool to subtract HIV-specific mortality from all-cause mortality after 1990. Smoothed over the trough in the 1990s and made the decline in under 5 mortality a little less severe. Extrapolated the trend between 2000 and 2017 to 2020."

# Reflect the scrolled viewport / selection recorded for this sheet at
# the time of the edit.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 97
$ws.Range("O207").Select()
